# Apply the "table Discount and Fee_Group(master-library) added" edit.
#
# This adds new data values to the existing "Fee details" (row 9) and
# "Discount" (row 10) table rows on Sheet1:
#   - Fee group (FEEGROUP)   column E
#   - Estimate time (20 mins) column G
#   - Status (completed)     column I
#   - Scholarship (SCHOLARSHIP) column E (Discount row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - "Fee details"
$ws.Range("E9").Value = "FEEGROUP"
$ws.Range("I9").Value = "completed"

# Row 10 - "Discount"
$ws.Range("E10").Value = "SCHOLARSHIP"

# Shared "20 mins" estimate-time values (added after the strings above so the
# shared-string table ends up in the same order as the source workbook).
$ws.Range("G9").Value = "20 mins"
$ws.Range("I10").Value = "completed"
$ws.Range("G10").Value = "20 mins"

# Update the view state to match: scrolled so E6 is the top-left visible
# cell, with G9 selected as the active cell.
$ws.Range("G9").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 5
